$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 251.75
$ws.Range("I58").Value = 200
$ws.Range("J58").Value = 303.5
$ws.Range("K58").Value = 600
$ws.Range("L58").Value = 910.5
$ws.Range("M58").Value = -450
$ws.Range("N58").Value = -1210.5

$ws.Range("H69").Value = 132503.25
$ws.Range("I69").Value = 5006.5
$ws.Range("J69").Value = 260000
$ws.Range("K69").Value = 15019.5
$ws.Range("L69").Value = 780000
$ws.Range("M69").Value = -14145.5
$ws.Range("N69").Value = -781748

$ws.Range("H72").Value = 132503.25
$ws.Range("I72").Value = 5006.5
$ws.Range("J72").Value = 260000
$ws.Range("K72").Value = 45058.5
$ws.Range("L72").Value = 2340000
$ws.Range("M72").Value = -40690.5
$ws.Range("N72").Value = -2348736

$ws.Range("H87").Value = 98853.5
$ws.Range("J87").Value = 98853.5
$ws.Range("L87").Value = 98853.5
$ws.Range("N87").Value = -101349.5

$ws.Range("H90").Value = 98853.5
$ws.Range("J90").Value = 98853.5
$ws.Range("L90").Value = 296560.5
$ws.Range("N90").Value = -309040.5

$ws.Range("H111").Value = 456.42856
$ws.Range("I111").Value = 399
$ws.Range("K111").Value = 1197
$ws.Range("M111").Value = 1870

$ws.Range("H112").Value = 2566.476
$ws.Range("I112").Value = 1930
$ws.Range("J112").Value = 2672.5557
$ws.Range("K112").Value = 5790
$ws.Range("L112").Value = 8017.6671
$ws.Range("M112").Value = -4682
$ws.Range("N112").Value = -10233.6671

$ws.Range("H132").Value = 3093.7
$ws.Range("I132").Value = 2392.25
$ws.Range("J132").Value = 5899.5
$ws.Range("K132").Value = 7176.75
$ws.Range("L132").Value = 17698.5
$ws.Range("M132").Value = -4646.75
$ws.Range("N132").Value = -22758.5

$ws.Range("H138").Value = 2968.625
$ws.Range("I138").Value = 2298.3333
$ws.Range("J138").Value = 4979.5
$ws.Range("K138").Value = 6894.999899999999
$ws.Range("L138").Value = 14938.5
$ws.Range("M138").Value = -1754.999899999999
$ws.Range("N138").Value = -25218.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4579.8
$ws.Range("I2").Value = 3999.5
$ws.Range("J2").Value = 4966.6665
$ws.Range("K2").Value = 3999.5
$ws.Range("L2").Value = 4966.6665
$ws.Range("M2").Value = -3886.5
$ws.Range("N2").Value = -5192.6665

$ws.Range("H63").Value = 12056.667
$ws.Range("I63").Value = 5126
$ws.Range("J63").Value = 17601.2
$ws.Range("K63").Value = 5126
$ws.Range("L63").Value = 17601.2
$ws.Range("M63").Value = -4440
$ws.Range("N63").Value = -18973.2

$ws.Range("H66").Value = 12056.667
$ws.Range("I66").Value = 5126
$ws.Range("J66").Value = 17601.2
$ws.Range("K66").Value = 25630
$ws.Range("L66").Value = 88006
$ws.Range("M66").Value = -22198
$ws.Range("N66").Value = -94870

$ws.Range("H88").Value = 2891.2
$ws.Range("I88").Value = 2125
$ws.Range("J88").Value = 3402
$ws.Range("K88").Value = 2125
$ws.Range("L88").Value = 3402
$ws.Range("M88").Value = -1719
$ws.Range("N88").Value = -4214

$ws.Range("H91").Value = 2891.2
$ws.Range("I91").Value = 2125
$ws.Range("J91").Value = 3402
$ws.Range("K91").Value = 2125
$ws.Range("L91").Value = 3402
$ws.Range("M91").Value = -721
$ws.Range("N91").Value = -6210

$ws.Range("H116").Value = 4579.8
$ws.Range("I116").Value = 3999.5
$ws.Range("J116").Value = 4966.6665
$ws.Range("K116").Value = 3999.5
$ws.Range("L116").Value = 4966.6665
$ws.Range("M116").Value = -1705.5
$ws.Range("N116").Value = -9554.666499999999

$ws.Range("H122").Value = 2179.8
$ws.Range("I122").Value = 2179.8
$ws.Range("K122").Value = 6539.400000000001
$ws.Range("M122").Value = -4089.400000000001

$ws.Range("H132").Value = 3785
$ws.Range("I132").Value = 3642
$ws.Range("K132").Value = 10926
$ws.Range("M132").Value = -8396

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4579.8
$ws.Range("I3").Value = 3999.5
$ws.Range("J3").Value = 4966.6665
$ws.Range("K3").Value = 3999.5
$ws.Range("L3").Value = 4966.6665
$ws.Range("M3").Value = -3885.5
$ws.Range("N3").Value = -5194.6665

$ws.Range("H99").Value = 1242.7142
$ws.Range("I99").Value = 1424.75
$ws.Range("K99").Value = 1424.75
$ws.Range("M99").Value = 73.25

$ws.Range("H105").Value = 3422.2
$ws.Range("I105").Value = 3422.2
$ws.Range("K105").Value = 3422.2
$ws.Range("M105").Value = -1675.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4081.3635
$ws.Range("I132").Value = 3679.6
$ws.Range("J132").Value = 4416.1665
$ws.Range("K132").Value = 11038.8
$ws.Range("L132").Value = 13248.4995
$ws.Range("M132").Value = -8508.799999999999
$ws.Range("N132").Value = -18308.4995

$ws.Range("H134").Value = 5779.5625
$ws.Range("I134").Value = 5359.8184
$ws.Range("K134").Value = 16079.4552
$ws.Range("M134").Value = -13544.4552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 22.166666
$ws.Range("I2").Value = 25.25
$ws.Range("J2").Value = 16
$ws.Range("K2").Value = 151.5
$ws.Range("L2").Value = 96
$ws.Range("M2").Value = -38.5
$ws.Range("N2").Value = -322

$ws.Range("H8").Value = 3003201.5
$ws.Range("I8").Value = 3003201.5
$ws.Range("K8").Value = 9009604.5
$ws.Range("M8").Value = -9009465.5

$ws.Range("H68").Value = 2963.6667
$ws.Range("J68").Value = 2963.6667
$ws.Range("L68").Value = 8891.000100000001
$ws.Range("N68").Value = -10513.0001

$ws.Range("H71").Value = 2963.6667
$ws.Range("J71").Value = 2963.6667
$ws.Range("L71").Value = 26673.0003
$ws.Range("N71").Value = -34785.0003

$ws.Range("H107").Value = 502
$ws.Range("J107").Value = 552.5
$ws.Range("L107").Value = 1657.5
$ws.Range("N107").Value = -5497.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H97").Value = 2024.9166
$ws.Range("I97").Value = 1449.8334
$ws.Range("J97").Value = 2600
$ws.Range("K97").Value = 1449.8334
$ws.Range("L97").Value = 2600
$ws.Range("M97").Value = -953.8334
$ws.Range("N97").Value = -3592

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7500
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5376

$ws.Range("H55").Value = 3600
$ws.Range("I55").Value = 4666.6665
$ws.Range("J55").Value = 2000
$ws.Range("K55").Value = 4666.6665
$ws.Range("L55").Value = 2000
$ws.Range("M55").Value = -4493.6665
$ws.Range("N55").Value = -2346

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 39992.332
$ws.Range("J74").Value = 37489
$ws.Range("L74").Value = 37489
$ws.Range("N74").Value = -39361

$ws.Range("H77").Value = 39992.332
$ws.Range("J77").Value = 37489
$ws.Range("L77").Value = 112467
$ws.Range("N77").Value = -121827

$ws.Range("H100").Value = 1526.9231
$ws.Range("I100").Value = 1257.125
$ws.Range("J100").Value = 1958.6
$ws.Range("K100").Value = 2514.25
$ws.Range("L100").Value = 3917.2
$ws.Range("M100").Value = -1973.25
$ws.Range("N100").Value = -4999.2

$ws.Range("H132").Value = 2561.125
$ws.Range("I132").Value = 1599
$ws.Range("J132").Value = 3138.4
$ws.Range("K132").Value = 4797
$ws.Range("L132").Value = 9415.200000000001
$ws.Range("M132").Value = -2267
$ws.Range("N132").Value = -14475.2
